$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Agosto de 2020 a las 09:07"

# Update per-country statistics (columns B-H) for the affected rows
# Row 4
$ws.Range("B4").Value = 4706059
$ws.Range("C4").Value = 170
$ws.Range("D4").Value = 2327791
$ws.Range("E4").Value = 2221516
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 156752

# Row 6
$ws.Range("B6").Value = 1701307
$ws.Range("C6").Value = 4253
$ws.Range("D6").Value = 1096893
$ws.Range("E6").Value = 567831
$ws.Range("G6").Value = 32
$ws.Range("H6").Value = 36583

# Row 53
$ws.Range("B53").Value = 38841
$ws.Range("C53").Value = 291
$ws.Range("D53").Value = 29557
$ws.Range("E53").Value = 8535
$ws.Range("G53").Value = 11
$ws.Range("H53").Value = 749

# Row 54
$ws.Range("B54").Value = 36710
$ws.Range("C54").Value = 35
$ws.Range("E54").Value = 9918
$ws.Range("G54").Value = 11
$ws.Range("H54").Value = 1283

# Row 72
$ws.Range("B72").Value = 17282
$ws.Range("C72").Value = 377
$ws.Range("D72").Value = 10202
$ws.Range("E72").Value = 6879

# Row 74
$ws.Range("D74").Value = 8396
$ws.Range("E74").Value = 7788

# Row 105
$ws.Range("B105").Value = 4526
$ws.Range("C105").Value = 21
$ws.Range("D105").Value = 3364
$ws.Range("E105").Value = 565
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 597

# Row 140
$ws.Range("B140").Value = 1238
$ws.Range("C140").Value = 7
$ws.Range("E140").Value = 154

# Row 143
$ws.Range("B143").Value = 1171
$ws.Range("C143").Value = 3
$ws.Range("D143").Value = 947
$ws.Range("E143").Value = 207

# Row 165
$ws.Range("B165").Value = 474
$ws.Range("C165").Value = 7
$ws.Range("E165").Value = 26
